$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

# Row 1, Col 1: "95 x 16" -> "52 x 54"
$t.Cell(1, 1).Range.Text = "52 x 54" + $nl + "  5    4" + $nl + "  ----" + $nl + "5|    |" + $nl + "2|    |"

# Row 1, Col 2: "56 x 59" -> "77 x 23"
$t.Cell(1, 2).Range.Text = "77 x 23" + $nl + "  2    3" + $nl + "  ----" + $nl + "7|    |" + $nl + "7|    |"

# Row 1, Col 3: "25 x 44" -> "67 x 15"
$t.Cell(1, 3).Range.Text = "67 x 15" + $nl + "  1    5" + $nl + "  ----" + $nl + "6|    |" + $nl + "7|    |"

# Row 2, Col 1: "15 x 56" -> "55 x 58"
$t.Cell(2, 1).Range.Text = "55 x 58" + $nl + "  5    8" + $nl + "  ----" + $nl + "5|    |" + $nl + "5|    |"

# Row 2, Col 2: "82 x 10" -> "14 x 32"
$t.Cell(2, 2).Range.Text = "14 x 32" + $nl + "  3    2" + $nl + "  ----" + $nl + "1|    |" + $nl + "4|    |"

# Row 2, Col 3: "17 x 24" -> "97 x 81"
$t.Cell(2, 3).Range.Text = "97 x 81" + $nl + "  8    1" + $nl + "  ----" + $nl + "9|    |" + $nl + "7|    |"

# Row 3, Col 1: "13 x 17" -> "55 x 97"
$t.Cell(3, 1).Range.Text = "55 x 97" + $nl + "  9    7" + $nl + "  ----" + $nl + "5|    |" + $nl + "5|    |"

# Row 3, Col 2: "65 x 79" -> "72 x 27"
$t.Cell(3, 2).Range.Text = "72 x 27" + $nl + "  2    7" + $nl + "  ----" + $nl + "7|    |" + $nl + "2|    |"

# Row 3, Col 3: "39 x 97" -> "98 x 33"
$t.Cell(3, 3).Range.Text = "98 x 33" + $nl + "  3    3" + $nl + "  ----" + $nl + "9|    |" + $nl + "8|    |"

# Row 4, Col 1: "97 x 20" -> "23 x 90"
$t.Cell(4, 1).Range.Text = "23 x 90" + $nl + "  9    0" + $nl + "  ----" + $nl + "2|    |" + $nl + "3|    |"

# Row 4, Col 2: "65 x 83" -> "49 x 53"
$t.Cell(4, 2).Range.Text = "49 x 53" + $nl + "  5    3" + $nl + "  ----" + $nl + "4|    |" + $nl + "9|    |"

# Row 4, Col 3: "66 x 99" -> "94 x 10"
$t.Cell(4, 3).Range.Text = "94 x 10" + $nl + "  1    0" + $nl + "  ----" + $nl + "9|    |" + $nl + "4|    |"

# Row 5, Col 1: "14 x 66" -> "65 x 17"
$t.Cell(5, 1).Range.Text = "65 x 17" + $nl + "  1    7" + $nl + "  ----" + $nl + "6|    |" + $nl + "5|    |"

# Row 5, Col 2: "60 x 93" -> "74 x 76"
$t.Cell(5, 2).Range.Text = "74 x 76" + $nl + "  7    6" + $nl + "  ----" + $nl + "7|    |" + $nl + "4|    |"

# Row 5, Col 3: "68 x 36" -> "60 x 58"
$t.Cell(5, 3).Range.Text = "60 x 58" + $nl + "  5    8" + $nl + "  ----" + $nl + "6|    |" + $nl + "0|    |"
